$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (147-151, 154, 340-345) ---
# Row 147
$ws.Cells.Item(147,2).Value = 5461474
$ws.Cells.Item(147,6).Value = "Legia Warsaw"
$ws.Cells.Item(147,7).Value = "Slask Wroclaw"
$ws.Cells.Item(147,8).Value = 3
$ws.Cells.Item(147,9).Value = 1
$ws.Cells.Item(147,11).Value = 1.7
$ws.Cells.Item(147,12).Value = 3.8
$ws.Cells.Item(147,13).Value = 4
$ws.Cells.Item(147,14).Value = 1.833
$ws.Cells.Item(147,15).Value = 3.8
$ws.Cells.Item(147,16).Value = 3.4
$ws.Cells.Item(147,17).Value = -0.5
$ws.Cells.Item(147,18).Value = 1.825
$ws.Cells.Item(147,19).Value = 2.025
$ws.Cells.Item(147,20).Value = 2.75
$ws.Cells.Item(147,21).Value = 1.9
$ws.Cells.Item(147,22).Value = 1.95
$ws.Cells.Item(147,23).Value = 0.833
$ws.Cells.Item(147,26).Value = 0.825
$ws.Cells.Item(147,28).Value = 0.8999999999999999
$ws.Cells.Item(147,29).Value = -1

# Row 148
$ws.Cells.Item(148,2).Value = 5428774
$ws.Cells.Item(148,6).Value = "Pogon Szczecin"
$ws.Cells.Item(148,7).Value = "Radomiak Radom"
$ws.Cells.Item(148,8).Value = 4
$ws.Cells.Item(148,9).Value = 0
$ws.Cells.Item(148,10).Value = "H"
$ws.Cells.Item(148,11).Value = 1.571
$ws.Cells.Item(148,12).Value = 4
$ws.Cells.Item(148,13).Value = 4.75
$ws.Cells.Item(148,14).Value = 1.533
$ws.Cells.Item(148,15).Value = 4.333
$ws.Cells.Item(148,16).Value = 4.75
$ws.Cells.Item(148,17).Value = -1
$ws.Cells.Item(148,18).Value = 1.875
$ws.Cells.Item(148,19).Value = 1.975
$ws.Cells.Item(148,20).Value = 3
$ws.Cells.Item(148,21).Value = 1.875
$ws.Cells.Item(148,22).Value = 1.975
$ws.Cells.Item(148,23).Value = 0.5329999999999999
$ws.Cells.Item(148,25).Value = -1
$ws.Cells.Item(148,26).Value = 0.875
$ws.Cells.Item(148,27).Value = -1
$ws.Cells.Item(148,28).Value = 0.875

# Row 149
$ws.Cells.Item(149,2).Value = 5467427
$ws.Cells.Item(149,6).Value = "Stal Mielec"
$ws.Cells.Item(149,7).Value = "Warta Poznan"
$ws.Cells.Item(149,8).Value = 1
$ws.Cells.Item(149,9).Value = 0
$ws.Cells.Item(149,11).Value = 2.375
$ws.Cells.Item(149,12).Value = 3.2
$ws.Cells.Item(149,13).Value = 2.8
$ws.Cells.Item(149,14).Value = 2.6
$ws.Cells.Item(149,15).Value = 3.1
$ws.Cells.Item(149,16).Value = 2.625
$ws.Cells.Item(149,17).Value = 0
$ws.Cells.Item(149,18).Value = 1.925
$ws.Cells.Item(149,19).Value = 1.925
$ws.Cells.Item(149,20).Value = 2.25
$ws.Cells.Item(149,21).Value = 1.975
$ws.Cells.Item(149,22).Value = 1.875
$ws.Cells.Item(149,23).Value = 1.6
$ws.Cells.Item(149,26).Value = 0.925
$ws.Cells.Item(149,28).Value = -1
$ws.Cells.Item(149,29).Value = 0.875

# Row 150
$ws.Cells.Item(150,2).Value = 5465446
$ws.Cells.Item(150,6).Value = "Cracovia Krakow"
$ws.Cells.Item(150,7).Value = "Wisla Plock"
$ws.Cells.Item(150,8).Value = 3
$ws.Cells.Item(150,10).Value = "H"
$ws.Cells.Item(150,11).Value = 2.15
$ws.Cells.Item(150,13).Value = 2.875
$ws.Cells.Item(150,14).Value = 2.25
$ws.Cells.Item(150,15).Value = 3.6
$ws.Cells.Item(150,16).Value = 2.7
$ws.Cells.Item(150,17).Value = -0.25
$ws.Cells.Item(150,18).Value = 2.05
$ws.Cells.Item(150,19).Value = 1.75
$ws.Cells.Item(150,20).Value = 2.5
$ws.Cells.Item(150,21).Value = 1.825
$ws.Cells.Item(150,22).Value = 2.025
$ws.Cells.Item(150,23).Value = 1.25
$ws.Cells.Item(150,24).Value = -1
$ws.Cells.Item(150,26).Value = 1.05
$ws.Cells.Item(150,27).Value = -1
$ws.Cells.Item(150,28).Value = 0.825
$ws.Cells.Item(150,29).Value = -1

# Row 151
$ws.Cells.Item(151,2).Value = 5461475
$ws.Cells.Item(151,6).Value = "Widzew Lodz"
$ws.Cells.Item(151,7).Value = "Korona Kielce"
$ws.Cells.Item(151,8).Value = 0
$ws.Cells.Item(151,9).Value = 3
$ws.Cells.Item(151,10).Value = "A"
$ws.Cells.Item(151,11).Value = 2.1
$ws.Cells.Item(151,12).Value = 3.3
$ws.Cells.Item(151,13).Value = 3.2
$ws.Cells.Item(151,14).Value = 2.375
$ws.Cells.Item(151,15).Value = 3.3
$ws.Cells.Item(151,16).Value = 2.7
$ws.Cells.Item(151,17).Value = 0
$ws.Cells.Item(151,18).Value = 1.8
$ws.Cells.Item(151,19).Value = 2.05
$ws.Cells.Item(151,20).Value = 2.5
$ws.Cells.Item(151,21).Value = 1.825
$ws.Cells.Item(151,22).Value = 2.025
$ws.Cells.Item(151,23).Value = -1
$ws.Cells.Item(151,25).Value = 1.7
$ws.Cells.Item(151,26).Value = -1
$ws.Cells.Item(151,27).Value = 1.05
$ws.Cells.Item(151,28).Value = 0.825

# Row 154
$ws.Cells.Item(154,2).Value = 5460884
$ws.Cells.Item(154,6).Value = "Miedz Legnica"
$ws.Cells.Item(154,7).Value = "Gornik Zabrze"
$ws.Cells.Item(154,8).Value = 0
$ws.Cells.Item(154,10).Value = "D"
$ws.Cells.Item(154,11).Value = 3.6
$ws.Cells.Item(154,13).Value = 1.909
$ws.Cells.Item(154,14).Value = 3
$ws.Cells.Item(154,15).Value = 3.5
$ws.Cells.Item(154,16).Value = 2.1
$ws.Cells.Item(154,17).Value = 0.25
$ws.Cells.Item(154,18).Value = 1.95
$ws.Cells.Item(154,19).Value = 1.9
$ws.Cells.Item(154,20).Value = 2.75
$ws.Cells.Item(154,21).Value = 1.975
$ws.Cells.Item(154,22).Value = 1.875
$ws.Cells.Item(154,23).Value = -1
$ws.Cells.Item(154,24).Value = 2.5
$ws.Cells.Item(154,26).Value = 0.475
$ws.Cells.Item(154,27).Value = -0.5
$ws.Cells.Item(154,28).Value = -1
$ws.Cells.Item(154,29).Value = 0.875

# Row 340
$ws.Cells.Item(340,8).Value = 3
$ws.Cells.Item(340,9).Value = 1
$ws.Cells.Item(340,10).Value = "H"
$ws.Cells.Item(340,14).Value = 2.05
$ws.Cells.Item(340,15).Value = 3.3
$ws.Cells.Item(340,17).Value = -0.5
$ws.Cells.Item(340,18).Value = 2
$ws.Cells.Item(340,19).Value = 1.85
$ws.Cells.Item(340,21).Value = 1.975
$ws.Cells.Item(340,22).Value = 1.875
$ws.Cells.Item(340,23).Value = 1.05
$ws.Cells.Item(340,24).Value = -1
$ws.Cells.Item(340,25).Value = -1
$ws.Cells.Item(340,26).Value = 1
$ws.Cells.Item(340,27).Value = -1
$ws.Cells.Item(340,28).Value = 0.9750000000000001
$ws.Cells.Item(340,29).Value = -1

# Row 341
$ws.Cells.Item(341,15).Value = 3
$ws.Cells.Item(341,21).Value = 1.95
$ws.Cells.Item(341,22).Value = 1.9

# Row 342
$ws.Cells.Item(342,21).Value = 1.975
$ws.Cells.Item(342,22).Value = 1.875

# Row 343
$ws.Cells.Item(343,18).Value = 2.05
$ws.Cells.Item(343,19).Value = 1.8
$ws.Cells.Item(343,21).Value = 1.9
$ws.Cells.Item(343,22).Value = 1.95

# Row 344
$ws.Cells.Item(344,18).Value = 1.8
$ws.Cells.Item(344,19).Value = 2.05
$ws.Cells.Item(344,21).Value = 1.875
$ws.Cells.Item(344,22).Value = 1.975

# Row 345
$ws.Cells.Item(345,2).Value = 6775549
$ws.Cells.Item(345,5).Value = 45346.66666666666
$ws.Cells.Item(345,6).Value = "Stal Mielec"
$ws.Cells.Item(345,7).Value = "Rakow Czestochowa"
$ws.Cells.Item(345,11).Value = 3.75
$ws.Cells.Item(345,12).Value = 3.5
$ws.Cells.Item(345,13).Value = 1.909
$ws.Cells.Item(345,14).Value = 4.5
$ws.Cells.Item(345,15).Value = 3.6
$ws.Cells.Item(345,16).Value = 1.75
$ws.Cells.Item(345,17).Value = 0.75
$ws.Cells.Item(345,18).Value = 1.875
$ws.Cells.Item(345,19).Value = 1.975
$ws.Cells.Item(345,21).Value = 1.875
$ws.Cells.Item(345,22).Value = 1.975

# --- Copy A/E column formatting (id bold+border style, Date numberformat) into the new rows ---
$srcA = $ws.Range("A345")
$srcE = $ws.Range("E345")
$srcA.Copy()
$ws.Range("A346").PasteSpecial(-4122)
$srcE.Copy()
$ws.Range("E346").PasteSpecial(-4122)
$srcA.Copy()
$ws.Range("A347").PasteSpecial(-4122)
$srcE.Copy()
$ws.Range("E347").PasteSpecial(-4122)
$srcA.Copy()
$ws.Range("A348").PasteSpecial(-4122)
$srcE.Copy()
$ws.Range("E348").PasteSpecial(-4122)
$srcA.Copy()
$ws.Range("A349").PasteSpecial(-4122)
$srcE.Copy()
$ws.Range("E349").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Add new rows (346-349) ---
# Row 346
$ws.Cells.Item(346,1).Value = 344
$ws.Cells.Item(346,2).Value = 6774881
$ws.Cells.Item(346,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(346,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(346,5).Value = 45347.35416666666
$ws.Cells.Item(346,6).Value = "MKS Puszcza Niepolomice"
$ws.Cells.Item(346,7).Value = "Zaglebie Lubin"
$ws.Cells.Item(346,11).Value = 3.3
$ws.Cells.Item(346,12).Value = 3.25
$ws.Cells.Item(346,13).Value = 2.15
$ws.Cells.Item(346,14).Value = 3.2
$ws.Cells.Item(346,15).Value = 3.25
$ws.Cells.Item(346,16).Value = 2.2
$ws.Cells.Item(346,17).Value = 0.25
$ws.Cells.Item(346,18).Value = 1.9
$ws.Cells.Item(346,19).Value = 1.95
$ws.Cells.Item(346,20).Value = 2.5
$ws.Cells.Item(346,21).Value = 1.975
$ws.Cells.Item(346,22).Value = 1.875
$ws.Cells.Item(346,23).Value = 0
$ws.Cells.Item(346,24).Value = 0
$ws.Cells.Item(346,25).Value = 0
$ws.Cells.Item(346,26).Value = 0
$ws.Cells.Item(346,27).Value = 0

# Row 347
$ws.Cells.Item(347,1).Value = 345
$ws.Cells.Item(347,2).Value = 6775551
$ws.Cells.Item(347,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(347,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(347,5).Value = 45347.45833333334
$ws.Cells.Item(347,6).Value = "Widzew Lodz"
$ws.Cells.Item(347,7).Value = "Gornik Zabrze"
$ws.Cells.Item(347,11).Value = 2.55
$ws.Cells.Item(347,12).Value = 3.2
$ws.Cells.Item(347,13).Value = 2.7
$ws.Cells.Item(347,14).Value = 2.3
$ws.Cells.Item(347,15).Value = 3.25
$ws.Cells.Item(347,16).Value = 3
$ws.Cells.Item(347,17).Value = -0.25
$ws.Cells.Item(347,18).Value = 2.025
$ws.Cells.Item(347,19).Value = 1.825
$ws.Cells.Item(347,20).Value = 2.5
$ws.Cells.Item(347,21).Value = 1.975
$ws.Cells.Item(347,22).Value = 1.875
$ws.Cells.Item(347,23).Value = 0
$ws.Cells.Item(347,24).Value = 0
$ws.Cells.Item(347,25).Value = 0
$ws.Cells.Item(347,26).Value = 0
$ws.Cells.Item(347,27).Value = 0

# Row 348
$ws.Cells.Item(348,1).Value = 346
$ws.Cells.Item(348,2).Value = 6775552
$ws.Cells.Item(348,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(348,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(348,5).Value = 45347.5625
$ws.Cells.Item(348,6).Value = "Korona Kielce"
$ws.Cells.Item(348,7).Value = "Legia Warsaw"
$ws.Cells.Item(348,11).Value = 3.5
$ws.Cells.Item(348,12).Value = 3.4
$ws.Cells.Item(348,13).Value = 2
$ws.Cells.Item(348,14).Value = 3.4
$ws.Cells.Item(348,15).Value = 3.4
$ws.Cells.Item(348,16).Value = 2.05
$ws.Cells.Item(348,17).Value = 0.25
$ws.Cells.Item(348,18).Value = 2.05
$ws.Cells.Item(348,19).Value = 1.8
$ws.Cells.Item(348,20).Value = 2.5
$ws.Cells.Item(348,21).Value = 2.025
$ws.Cells.Item(348,22).Value = 1.825
$ws.Cells.Item(348,23).Value = 0
$ws.Cells.Item(348,24).Value = 0
$ws.Cells.Item(348,25).Value = 0
$ws.Cells.Item(348,26).Value = 0
$ws.Cells.Item(348,27).Value = 0

# Row 349
$ws.Cells.Item(349,1).Value = 347
$ws.Cells.Item(349,2).Value = 6775550
$ws.Cells.Item(349,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(349,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(349,5).Value = 45348.625
$ws.Cells.Item(349,6).Value = "Warta Poznan"
$ws.Cells.Item(349,7).Value = "Radomiak Radom"
$ws.Cells.Item(349,11).Value = 2.5
$ws.Cells.Item(349,12).Value = 3.25
$ws.Cells.Item(349,13).Value = 2.75
$ws.Cells.Item(349,14).Value = 2.45
$ws.Cells.Item(349,15).Value = 3.3
$ws.Cells.Item(349,16).Value = 2.875
$ws.Cells.Item(349,17).Value = 0
$ws.Cells.Item(349,18).Value = 1.775
$ws.Cells.Item(349,19).Value = 2.1
$ws.Cells.Item(349,20).Value = 2.25
$ws.Cells.Item(349,21).Value = 1.875
$ws.Cells.Item(349,22).Value = 1.975
$ws.Cells.Item(349,23).Value = 0
$ws.Cells.Item(349,24).Value = 0
$ws.Cells.Item(349,25).Value = 0
$ws.Cells.Item(349,26).Value = 0
$ws.Cells.Item(349,27).Value = 0
